$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '25.986.67'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +0.12%  '

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.642.09'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -0.02%  '

$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.44%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '215.01'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.16%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.5064'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -0.42%  '

$ws.Range('E7').Value = '  -0.30%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2579'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.58%  '

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06364'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.44%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '19.90'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +1.78%  '

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07739'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.51%  '

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '4.301'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.08%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.631.58'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -0.66%  '

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.5475'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +0.44%  '

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.0₅7757'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -1.15%  '

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '64.30'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.74%  '

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '26.012.62'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.11%  '

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.22%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '4.475'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.92%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '196.69'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.63%  '

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '9.985'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.28%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.159'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.92%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.39%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '1.894'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +0.93%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '142.35'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +1.03%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.1263'
$c.Style = 'Normal'

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '6.875'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.15%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '15.66'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.46%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.241'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +0.34%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.04897'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -2.61%  '

$ws.Range('E31').Value = '  +0.28%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.210'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +0.74%  '

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.553'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.77%  '

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '2.377'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +0.71%  '

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.9191'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +2.71%  '

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.567'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -0.80%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.5552'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +0.98%  '

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '1.132.89'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +0.51%  '

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.01570'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +1.06%  '

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.30%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '5.609'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.33%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.8037'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -1.63%  '

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '98.72'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -1.17%  '

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.779.11'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.05%  '

$ws.Range('E45').Value = '  -9.74%  '

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.4526'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.14%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '55.35'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +1.01%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.27%  '

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.05187'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +2.12%  '

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '7.602'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +2.71%  '

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -0.03%  '

